$wb = $excel.ActiveWorkbook

# NOTE 1: "Vector_bf" and "Vector_BF" differ only by case. Looking sheets up
# by name (Worksheets.Item("...")) gets confused once both of those names
# are used in the same script, so every sheet below is addressed by its
# 1-based index instead of its name:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha
#
# NOTE 2: every value in this workbook is stored as *text* (shared string),
# even the ones that look numeric (e.g. "-2.9", "0.97", "0"). Assigning a
# numeric-looking string straight to Range.Value auto-converts it to a real
# number, so for those cells we temporarily switch the cell to text format
# ("@"), assign the value, then reset the cell style back to "Normal" so no
# stray number formatting is left behind on the cell.

# Restricciones_del_lider
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "1.9 - x"

$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "-2.9"
$ws2.Range("B2").Style = "Normal"

$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "0.83"
$ws2.Range("D2").Style = "Normal"

$ws2.Range("A3").Value = "-1.9 + x"

$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "0.8999999999999999"
$ws2.Range("B3").Style = "Normal"

$ws2.Range("D3").NumberFormat = "@"
$ws2.Range("D3").Value = "0.08"
$ws2.Range("D3").Style = "Normal"

# Restricciones_del_follower
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "0.5062500000000001 - 3.375y"

$ws3.Range("B2").NumberFormat = "@"
$ws3.Range("B2").Value = "-1.50625"
$ws3.Range("B2").Style = "Normal"

$ws3.Range("D2").NumberFormat = "@"
$ws3.Range("D2").Value = "0.01"
$ws3.Range("D2").Style = "Normal"

$ws3.Range("E2").NumberFormat = "@"
$ws3.Range("E2").Value = "-4.6000000000000005"
$ws3.Range("E2").Style = "Normal"

$ws3.Range("F2").NumberFormat = "@"
$ws3.Range("F2").Value = "-8.100000000000001"
$ws3.Range("F2").Style = "Normal"

$ws3.Range("A3").NumberFormat = "@"
$ws3.Range("A3").Value = "0"
$ws3.Range("A3").Style = "Normal"

$ws3.Range("B3").NumberFormat = "@"
$ws3.Range("B3").Value = "-1"
$ws3.Range("B3").Style = "Normal"

$ws3.Range("D3").NumberFormat = "@"
$ws3.Range("D3").Value = "0.97"
$ws3.Range("D3").Style = "Normal"

$ws3.Range("E3").NumberFormat = "@"
$ws3.Range("E3").Value = "0"
$ws3.Range("E3").Style = "Normal"

$ws3.Range("F3").NumberFormat = "@"
$ws3.Range("F3").Value = "0"
$ws3.Range("F3").Style = "Normal"

# Punto_modificado
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("A2").NumberFormat = "@"
$ws4.Range("A2").Value = "1.9"
$ws4.Range("A2").Style = "Normal"

$ws4.Range("B2").NumberFormat = "@"
$ws4.Range("B2").Value = "0.15"
$ws4.Range("B2").Style = "Normal"

# Vector_bf
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("A2").NumberFormat = "@"
$ws5.Range("A2").Value = "0.597"
$ws5.Range("A2").Style = "Normal"

# Vector_BF
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("A2").NumberFormat = "@"
$ws6.Range("A2").Value = "2.3200000000000003"
$ws6.Range("A2").Style = "Normal"

$ws6.Range("A3").NumberFormat = "@"
$ws6.Range("A3").Value = "-11.629000000000001"
$ws6.Range("A3").Style = "Normal"

# Vector_Alpha
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("A2").Value = 2.4000000000000004
